# Generate Report for Handback
# Rows 4 and 5 (6daae5bf-... and b1694271-...) move from "Ready for handoff"
# to "Handed back: in sync with en-US" on all three sheets, and the
# per-language handback info (Latest Target File / Latest Handback File /
# Latest Handback DateTime) gets populated for those two rows.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview": columns E (zh-cn) and F (de-de) show the status.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $statusHandedBack
$wsOverview.Range("F4").Value = $statusHandedBack
$wsOverview.Range("E5").Value = $statusHandedBack
$wsOverview.Range("F5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $statusHandedBack
$wsZhCn.Range("I4").Value = "6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md"
$wsZhCn.Range("J4").Value = "6daae5bf-6b02-45e6-9a1c-31aacda4d54b.33f18f7d6bb71bf02c5e79d7b35d2694f0083ec3.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-20 18:38:33"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/792213980b0bfbe8e542f58a6bfadaf166aa2dfc/e2e/6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md", "", "", "6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md")

$wsZhCn.Range("C5").Value = $statusHandedBack
$wsZhCn.Range("I5").Value = "b1694271-d532-4fbc-b30f-dfcb4679f371.md"
$wsZhCn.Range("J5").Value = "b1694271-d532-4fbc-b30f-dfcb4679f371.0d33f40b4aef128c47179ab469e4b40af6c0bd81.zh-cn.xlf"
$wsZhCn.Range("K5").Value = "2016-08-20 18:38:33"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/792213980b0bfbe8e542f58a6bfadaf166aa2dfc/e2e/b1694271-d532-4fbc-b30f-dfcb4679f371.md", "", "", "b1694271-d532-4fbc-b30f-dfcb4679f371.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $statusHandedBack
$wsDeDe.Range("I4").Value = "6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md"
$wsDeDe.Range("J4").Value = "6daae5bf-6b02-45e6-9a1c-31aacda4d54b.33f18f7d6bb71bf02c5e79d7b35d2694f0083ec3.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-20 18:38:39"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8bc580f91c2d30133fd64349da9c2156479c8ac9/e2e/6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md", "", "", "6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md")

$wsDeDe.Range("C5").Value = $statusHandedBack
$wsDeDe.Range("I5").Value = "b1694271-d532-4fbc-b30f-dfcb4679f371.md"
$wsDeDe.Range("J5").Value = "b1694271-d532-4fbc-b30f-dfcb4679f371.0d33f40b4aef128c47179ab469e4b40af6c0bd81.de-de.xlf"
$wsDeDe.Range("K5").Value = "2016-08-20 18:38:39"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8bc580f91c2d30133fd64349da9c2156479c8ac9/e2e/b1694271-d532-4fbc-b30f-dfcb4679f371.md", "", "", "b1694271-d532-4fbc-b30f-dfcb4679f371.md")

Write-Output "Done applying handback updates"
